$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.954.73"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.225.23"
$ws.Range("E3").Value = "  -1.31%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.11"
$ws.Range("E5").Value = "  -1.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("E6").Value = "  +0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.39"
$ws.Range("E7").Value = "  -4.77%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.600"
$ws.Range("E9").Value = "  -4.43%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.40"
$ws.Range("E10").Value = "  -7.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0956"
$ws.Range("E11").Value = "  +0.34%  "

$ws.Range("E12").Value = "  +0.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.93"
$ws.Range("E13").Value = "  -5.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.554.61"
$ws.Range("E14").Value = "  -1.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.33"
$ws.Range("E15").Value = "  -2.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("E16").Value = "  -3.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.245.47"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.782.11"
$ws.Range("E18").Value = "  -1.49%  "

$ws.Range("E19").Value = "  +3.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.26"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.84"
$ws.Range("E21").Value = "  +0.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.29"
$ws.Range("E22").Value = "  +16.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.26"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  -8.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.39"
$ws.Range("E26").Value = "  -3.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.63"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.29"
$ws.Range("E28").Value = "  -1.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.05"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.53"
$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("E32").Value = "  +3.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0799"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.09"
$ws.Range("E34").Value = "  -5.03%  "

$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.108"
$ws.Range("E36").Value = "  -9.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.30"
$ws.Range("E37").Value = "  -7.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0303"
$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.45"
$ws.Range("E39").Value = "  -8.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.13"
$ws.Range("E40").Value = "  -3.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.98"
$ws.Range("E41").Value = "  +0.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.66"
$ws.Range("E42").Value = "  -2.95%  "

$ws.Range("E43").Value = "  -2.93%  "

$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.94"
$ws.Range("E45").Value = "  -3.10%  "

$ws.Range("E46").Value = "  -3.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.34"
$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("E49").Value = "  -2.64%  "

$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.426.60"
$ws.Range("E51").Value = "  -1.72%  "
